# The deck ships with two embedded themes: the slide master currently uses
# the "Integral" colour theme while the notes master carries the stock
# "Office Theme" colours. This edit swaps the colour scheme that is live on
# the design (slide master) over to the standard "Office Theme" palette.
#
# Helper: turn an "RRGGBB" hex string into the BGR-packed decimal value the
# PowerPoint ColorScheme/RGB property expects.
function Hex-ToRgbValue([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

$p = $ppt.ActivePresentation

# Target palette: the standard Office Theme 12-colour scheme, in the
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink order used by
# ColorScheme.Colors()/Item().
$officeTheme = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$cs = $p.SlideMaster.ColorScheme
for ($i = 1; $i -le $officeTheme.Count; $i++) {
    $cs.Colors($i).RGB = Hex-ToRgbValue $officeTheme[$i - 1]
}
